$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3202414.5
$ws.Range("J17").Value = 3202414.5
$ws.Range("L17").Value = 9607243.5
$ws.Range("N17").Value = -9607579.5

$ws.Range("H98").Value = 934.5833
$ws.Range("I98").Value = 837.7273
$ws.Range("K98").Value = 837.7273
$ws.Range("M98").Value = 660.2727

$ws.Range("H112").Value = 13515971
$ws.Range("J112").Value = 14288016
$ws.Range("L112").Value = 42864048
$ws.Range("N112").Value = -42866264

$ws.Range("H115").Value = 485
$ws.Range("I115").Value = 485
$ws.Range("K115").Value = 1455
$ws.Range("M115").Value = 112

$ws.Range("H122").Value = 934.5833
$ws.Range("I122").Value = 837.7273
$ws.Range("K122").Value = 2513.1819
$ws.Range("M122").Value = -63.18190000000004

$ws.Range("H136").Value = 90000
$ws.Range("J136").Value = 90000
$ws.Range("L136").Value = 90000
$ws.Range("N136").Value = -100200

$ws.Range("H138").Value = 23261406
$ws.Range("I138").Value = 1668.1111
$ws.Range("K138").Value = 5004.3333
$ws.Range("M138").Value = 135.6666999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 84188
$ws.Range("I31").Value = 14561.857
$ws.Range("K31").Value = 14561.857
$ws.Range("M31").Value = -14267.857

$ws.Range("H32").Value = 4055.6985
$ws.Range("I32").Value = 4056.5967
$ws.Range("K32").Value = 4056.5967
$ws.Range("M32").Value = -3769.5967

$ws.Range("H38").Value = 3000
$ws.Range("I38").Value = 3000
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 3000
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -2533
$ws.Range("N38").ClearContents()

$ws.Range("H97").Value = 516.13794
$ws.Range("I97").Value = 390.75
$ws.Range("J97").Value = 1118
$ws.Range("K97").Value = 390.75
$ws.Range("L97").Value = 1118
$ws.Range("M97").Value = 105.25
$ws.Range("N97").Value = -2110

$ws.Range("H110").Value = 28615.076
$ws.Range("I110").Value = 32908.727
$ws.Range("K110").Value = 32908.727
$ws.Range("M110").Value = -30863.727

$ws.Range("H112").Value = 18385.75
$ws.Range("J112").Value = 18385.75
$ws.Range("L112").Value = 18385.75
$ws.Range("N112").Value = -21339.75

$ws.Range("H138").Value = 199999
$ws.Range("J138").Value = 199999
$ws.Range("L138").Value = 199999
$ws.Range("N138").Value = -210279

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 35252.668
$ws.Range("I82").Value = 15340.333
$ws.Range("J82").Value = 75077.336
$ws.Range("K82").Value = 15340.333
$ws.Range("L82").Value = 75077.336
$ws.Range("M82").Value = -14957.333
$ws.Range("N82").Value = -75843.336

$ws.Range("H85").Value = 35252.668
$ws.Range("I85").Value = 15340.333
$ws.Range("J85").Value = 75077.336
$ws.Range("K85").Value = 15340.333
$ws.Range("L85").Value = 75077.336
$ws.Range("M85").Value = -14014.333
$ws.Range("N85").Value = -77729.336

$ws.Range("H99").Value = 3059.75
$ws.Range("I99").Value = 1633.6
$ws.Range("K99").Value = 1633.6
$ws.Range("M99").Value = -135.5999999999999

$ws.Range("H137").Value = 118284.86
$ws.Range("I137").Value = 114497
$ws.Range("K137").Value = 114497
$ws.Range("M137").Value = -109397

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws.Range("H31").Value = 4435.246
$ws.Range("I31").Value = 2917.5881
$ws.Range("J31").Value = 6346.3706
$ws.Range("K31").Value = 2917.5881
$ws.Range("L31").Value = 6346.3706
$ws.Range("M31").Value = -2622.5881
$ws.Range("N31").Value = -6936.3706

$ws.Range("H34").Value = 4435.246
$ws.Range("I34").Value = 2917.5881
$ws.Range("J34").Value = 6346.3706
$ws.Range("K34").Value = 2917.5881
$ws.Range("L34").Value = 6346.3706
$ws.Range("M34").Value = -2715.5881
$ws.Range("N34").Value = -6750.3706

$ws.Range("H132").Value = 44849.02
$ws.Range("I132").Value = 63902.305
$ws.Range("K132").Value = 191706.915
$ws.Range("M132").Value = -189176.915

$ws.Range("H134").Value = 2308.9048
$ws.Range("I134").Value = 2121.2974
$ws.Range("K134").Value = 6363.8922
$ws.Range("M134").Value = -3828.8922

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1476.1666
$ws.Range("I18").Value = 310.25
$ws.Range("J18").Value = 3808
$ws.Range("K18").Value = 930.75
$ws.Range("L18").Value = 11424
$ws.Range("M18").Value = -761.75
$ws.Range("N18").Value = -11762

$ws.Range("H130").Value = 2288
$ws.Range("I130").Value = 1183.3334
$ws.Range("K130").Value = 3550.0002
$ws.Range("M130").Value = 1469.9998

$ws.Range("H138").Value = 2358.25
$ws.Range("J138").Value = 2922
$ws.Range("L138").Value = 8766
$ws.Range("N138").Value = -19046

$ws.Range("H139").Value = 2715.2144
$ws.Range("I139").Value = 2662.25
$ws.Range("K139").Value = 7986.75
$ws.Range("M139").Value = -2846.75

$ws.Range("H141").Value = 13741.375
$ws.Range("I141").Value = 15299.667
$ws.Range("K141").Value = 45899.001
$ws.Range("M141").Value = -40719.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1667.7142
$ws.Range("I80").Value = 1667.7142
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1667.7142
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -669.7141999999999
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 1667.7142
$ws.Range("I83").Value = 1667.7142
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 8338.571
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -3346.571
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5879.8184
$ws.Range("I7").Value = 6128
$ws.Range("J7").Value = 5786.75
$ws.Range("K7").Value = 6128
$ws.Range("L7").Value = 5786.75
$ws.Range("M7").Value = -6016
$ws.Range("N7").Value = -6010.75

$ws.Range("H22").Value = 1617.8
$ws.Range("I22").Value = 1272.25
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 1272.25
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -977.25
$ws.Range("N22").Value = -3590

$ws.Range("H27").Value = 1617.8
$ws.Range("I27").Value = 1272.25
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 1272.25
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -1165.25
$ws.Range("N27").Value = -3214

$ws.Range("H68").Value = 2366.6667
$ws.Range("I68").Value = 2350
$ws.Range("K68").Value = 2350
$ws.Range("M68").Value = -1601

$ws.Range("H71").Value = 2366.6667
$ws.Range("I71").Value = 2350
$ws.Range("K71").Value = 11750
$ws.Range("M71").Value = -8006

$ws.Range("H82").Value = 4109.5
$ws.Range("I82").Value = 2344
$ws.Range("K82").Value = 2344
$ws.Range("M82").Value = -1983

$ws.Range("H85").Value = 4109.5
$ws.Range("I85").Value = 2344
$ws.Range("K85").Value = 2344
$ws.Range("M85").Value = -1096

$ws.Range("H126").Value = 5879.8184
$ws.Range("I126").Value = 6128
$ws.Range("J126").Value = 5786.75
$ws.Range("K126").Value = 18384
$ws.Range("L126").Value = 17360.25
$ws.Range("M126").Value = -15914
$ws.Range("N126").Value = -22300.25

$ws.Range("H132").Value = 7806.724
$ws.Range("I132").Value = 3550.8572
$ws.Range("K132").Value = 10652.5716
$ws.Range("M132").Value = -8122.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1216.9395
$ws.Range("I81").Value = 1364
$ws.Range("J81").Value = 1078.5294
$ws.Range("K81").Value = 2728
$ws.Range("L81").Value = 2157.0588
$ws.Range("M81").Value = -1667
$ws.Range("N81").Value = -4279.0588

$ws.Range("H84").Value = 1216.9395
$ws.Range("I84").Value = 1364
$ws.Range("J84").Value = 1078.5294
$ws.Range("K84").Value = 13640
$ws.Range("L84").Value = 10785.294
$ws.Range("M84").Value = -8336
$ws.Range("N84").Value = -21393.294

$ws.Range("H132").Value = 1577.3182
$ws.Range("I132").Value = 1097.2667
$ws.Range("J132").Value = 2606
$ws.Range("K132").Value = 3291.800099999999
$ws.Range("L132").Value = 7818
$ws.Range("M132").Value = -761.8000999999995
$ws.Range("N132").Value = -12878
